# Add 2022-Q3 data:
#  - Insert a new worksheet "2022-Q3" right before the "2022-Q2" sheet.
#  - Populate it with the fund-holdings table for that quarter.
#  - Update the "总计" (totals) sheet: insert a new row for 2022-Q3 at the
#    top of the data (row 2), pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet before "2022-Q2"
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$newWs = $wb.Worksheets.Add($q2Sheet)
$newWs.Name = "2022-Q3"

# Header row
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

$newWs.Range("B1:H1").Font.Bold = $true
$newWs.Range("B1:H1").HorizontalAlignment = -4108
$newWs.Range("B1:H1").VerticalAlignment = -4160
$newWs.Range("B1:H1").Borders.LineStyle = 1

# Columns B..G hold text data (fund codes must keep leading zeros, and the
# numeric-looking figures are stored as text in the source data)
$newWs.Range("B2:G6").NumberFormat = "@"

$data = @(
  @("010490", "鹏华高质量增长混合A",       "12.74", "93.98", "7.97", "1.0154", 4),
  @("009023", "鹏华稳健回报混合",           "4.12",  "94.39", "8.15", "0.3358", 4),
  @("007731", "民生加银持续成长混合A",     "3.22",  "94.57", "5.42", "0.1745", 10),
  @("007732", "民生加银持续成长混合C",     "1.89",  "94.57", "5.42", "0.1024", 10),
  @("010491", "鹏华高质量增长混合C",       "0.44",  "93.98", "7.97", "0.0351", 4)
)

$r = 2
foreach ($row in $data) {
    $newWs.Range("A$r").Value = ($r - 2)
    $newWs.Range("B$r").Value = $row[0]
    $newWs.Range("C$r").Value = $row[1]
    $newWs.Range("D$r").Value = $row[2]
    $newWs.Range("E$r").Value = $row[3]
    $newWs.Range("F$r").Value = $row[4]
    $newWs.Range("G$r").Value = $row[5]
    $newWs.Range("H$r").Value = $row[6]
    $r++
}

$newWs.Range("A2:A6").Font.Bold = $true
$newWs.Range("A2:A6").HorizontalAlignment = -4108
$newWs.Range("A2:A6").VerticalAlignment = -4160
$newWs.Range("A2:A6").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: shift existing rows down one and
#    insert the 2022-Q3 totals as the new row 2.
# ---------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")

# Work from the bottom up so we never clobber data before reading it.
$totalWs.Range("A7").Value = 5
$totalWs.Range("B7").Value = "2021-Q2"
$totalWs.Range("C7").Value = 2
$totalWs.Range("D7").Value = 0.02
$totalWs.Range("A7").Font.Bold = $true
$totalWs.Range("A7").HorizontalAlignment = -4108
$totalWs.Range("A7").VerticalAlignment = -4160
$totalWs.Range("A7").Borders.LineStyle = 1

$totalWs.Range("A6").Value = 4
$totalWs.Range("B6").Value = "2021-Q3"
$totalWs.Range("C6").Value = 6
$totalWs.Range("D6").Value = 1.3

$totalWs.Range("A5").Value = 3
$totalWs.Range("B5").Value = "2021-Q4"
$totalWs.Range("C5").Value = 12
$totalWs.Range("D5").Value = 2.76

$totalWs.Range("A4").Value = 2
$totalWs.Range("B4").Value = "2022-Q1"
$totalWs.Range("C4").Value = 3
$totalWs.Range("D4").Value = 0.74

$totalWs.Range("A3").Value = 1
$totalWs.Range("B3").Value = "2022-Q2"
$totalWs.Range("C3").Value = 6
$totalWs.Range("D3").Value = 1.93

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q3"
$totalWs.Range("C2").Value = 5
$totalWs.Range("D2").Value = 1.66
